# Insert a new data row at row 239 (above the current row 239), shifting all
# subsequent rows (239-295) down by one (to 240-296). This mirrors a new
# weekly price observation being added to the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 239; existing row 239 (and everything
# below it) shifts down to row 240, etc. Formatting (e.g. the date style on
# column D) is inherited from the row being pushed down, matching the
# original workbook's style usage.
$ws.Rows("239:239").Insert()

# Populate the newly inserted row 239 with the observation. Most fields are
# identical to the record that used to occupy row 239 (now row 240); only
# the date (D), quality (I) and volume (J) differ for this new entry.
$ws.Cells.Item(239, 1).Value2 = 4
$ws.Cells.Item(239, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(239, 3).Value2 = "Los Lagos"
$ws.Cells.Item(239, 4).Value2 = 44782
$ws.Cells.Item(239, 5).Value2 = 10
$ws.Cells.Item(239, 6).Value2 = 100112017
$ws.Cells.Item(239, 7).Value2 = "Apio"
$ws.Cells.Item(239, 8).Value2 = "Americana (o)"
$ws.Cells.Item(239, 9).Value2 = "Segunda"
$ws.Cells.Item(239, 10).Value2 = 25
$ws.Cells.Item(239, 11).Value2 = 12000
$ws.Cells.Item(239, 12).Value2 = 12000
$ws.Cells.Item(239, 13).Value2 = 12000
$ws.Cells.Item(239, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(239, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(239, 16).Value2 = 2000
$ws.Cells.Item(239, 17).Value2 = 6
$ws.Cells.Item(239, 18).Value2 = "Hortaliza"
